$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp in the title cell (A1) ---
$title = $ws.Cells.Item(1, 1).Value()
$ws.Cells.Item(1, 1).Value = $title.Replace("06:16", "06:46")

# --- Update India's row (row 43) counts ---
$ws.Cells.Item(43, 5).Value = 512   # Casos activos (E)
$ws.Cells.Item(43, 7).Value = 0     # Casos criticos (G)
$ws.Cells.Item(43, 8).Value = 10    # Muertes (H)

# --- Countries table: Kazajistan moves up (new data) right after Bielorrusia,
#     pushing Afganistan / Guadalupe / Costa de Marfil down one row each,
#     keeping their existing data unchanged. Georgia (row 105) and below
#     are untouched. ---

# Row 104 (was Kazajistan) becomes Costa de Marfil's old data
$ws.Cells.Item(104, 1).Value = "Costa de Marfil"
$ws.Cells.Item(104, 2).Value = 73
$ws.Cells.Item(104, 3).Value = 0
$ws.Cells.Item(104, 4).Value = 2
$ws.Cells.Item(104, 5).Value = 71
$ws.Cells.Item(104, 6).Value = 0
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = 0

# Row 103 (was Costa de Marfil) becomes Guadalupe's old data
$ws.Cells.Item(103, 1).Value = "Guadalupe"
$ws.Cells.Item(103, 2).Value = 73
$ws.Cells.Item(103, 3).Value = 0
$ws.Cells.Item(103, 4).Value = 0
$ws.Cells.Item(103, 5).Value = 72
$ws.Cells.Item(103, 6).Value = 4
$ws.Cells.Item(103, 7).Value = 0
$ws.Cells.Item(103, 8).Value = 1

# Row 102 (was Guadalupe) becomes Afganistan's old data
$ws.Cells.Item(102, 1).Value = "Afganistan"
$ws.Cells.Item(102, 2).Value = 74
$ws.Cells.Item(102, 3).Value = 0
$ws.Cells.Item(102, 4).Value = 1
$ws.Cells.Item(102, 5).Value = 72
$ws.Cells.Item(102, 6).Value = 0
$ws.Cells.Item(102, 7).Value = 0
$ws.Cells.Item(102, 8).Value = 1

# Row 101 (was Afganistan) becomes the new Kazajistan row with updated data
$ws.Cells.Item(101, 1).Value = "Kazajistan"
$ws.Cells.Item(101, 2).Value = 79
$ws.Cells.Item(101, 3).Value = 7
$ws.Cells.Item(101, 4).Value = 0
$ws.Cells.Item(101, 5).Value = 79
$ws.Cells.Item(101, 6).Value = 0
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 0
